$d = $word.ActiveDocument

$replacements = @(
    @{old="453×9=4077"; new="130×2=260"},
    @{old="448×6=2688"; new="978×6=5868"},
    @{old="964×5=4820"; new="282×8=2256"},
    @{old="579×6=3474"; new="501×5=2505"},
    @{old="746×2=1492"; new="340×3=1020"},
    @{old="640×7=4480"; new="429×5=2145"},
    @{old="729×4=2916"; new="854×7=5978"},
    @{old="807×3=2421"; new="988×2=1976"},
    @{old="958×8=7664"; new="884×6=5304"},
    @{old="748×9=6732"; new="676×5=3380"},
    @{old="598×9=5382"; new="952×7=6664"},
    @{old="762×9=6858"; new="420×9=3780"},
    @{old="546×9=4914"; new="959×5=4795"},
    @{old="105×7=735";  new="764×7=5348"},
    @{old="938×2=1876"; new="414×5=2070"},
    @{old="788×7=5516"; new="574×9=5166"},
    @{old="519×4=2076"; new="527×2=1054"},
    @{old="476×7=3332"; new="599×8=4792"},
    @{old="323×4=1292"; new="526×2=1052"},
    @{old="805×3=2415"; new="852×7=5964"},
    @{old="439×3=1317"; new="746×6=4476"},
    @{old="321×2=642";  new="509×9=4581"},
    @{old="137×9=1233"; new="705×9=6345"},
    @{old="693×6=4158"; new="829×2=1658"},
    @{old="401×8=3208"; new="151×7=1057"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $r.new, 2)
}
